# Horarios actualizados Línea 141 - 113
# Applies the scrape update from 20:45:46 to 20:52:24 to the "LP1912" workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

$oldStamp = "Última actualización: 20:45:46"
$newStamp = "Última actualización: 20:52:24"

# --- Sheet "LP1912" ---

# Update header info
$ws1.Range("A2").Value = $newStamp
$ws1.Range("A3").Value = "Total filas: 153"

# Rows 139 and 140 got swapped in the new scrape
$ws1.Range("A139").Value = "20:45:46"
$ws1.Range("B139").Value = "21:06"
$ws1.Range("C139").Value = "16_SANTA ANA"
$ws1.Range("D139").Value = 21
$ws1.Range("E139").Value = "LP1912"

$ws1.Range("A140").Value = "19:47:58"
$ws1.Range("B140").Value = "21:06"
$ws1.Range("C140").Value = "27_EL RETIRO"
$ws1.Range("D140").Value = 79
$ws1.Range("E140").Value = "LP1912"

# New row 158 appended with the latest scraped data
$ws1.Range("A158").Value = "20:52:24"
$ws1.Range("B158").Value = "22:49"
$ws1.Range("C158").Value = "14_ABASTO"
$ws1.Range("D158").Value = 117
$ws1.Range("E158").Value = "LP1912"

# --- Sheet "LP1912-215" ---
$ws2.Range("A2").Value = $newStamp

# --- Sheet "6203-6173" ---
$ws3.Range("A2").Value = $newStamp
